# The presentation originally contains 7 slides. The two "Tutoriale" slides
# covering PostgreSQL installation/SQL tutorials (old slides 3 and 4) are no
# longer wanted and must be removed, leaving the R/RStudio tutorial slides
# (old slides 5, 6, 7) shifted up to become slides 3, 4 and 5.
#
# Slides are located by their title text rather than a hard-coded index so
# the script is robust even if slide ordering assumptions differ slightly.

$p = $ppt.ActivePresentation

function Get-SlideTitleText($slide) {
    try {
        if ($slide.Shapes.Count -ge 1) {
            return $slide.Shapes.Item(1).TextFrame.TextRange.Text
        }
    } catch {
    }
    return ""
}

# Collect the indexes (1-based) of slides whose title indicates they are the
# PostgreSQL-installation / SQL-in-PostgreSQL tutorial slides that must be
# deleted.
$toDelete = New-Object System.Collections.ArrayList

for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $title = Get-SlideTitleText $p.Slides.Item($i)
    if ($title -match "instalare/lucru" -or $title -match "SQL") {
        [void]$toDelete.Add($i)
    }
}

# Delete from the highest index down to the lowest so earlier indexes stay
# valid while we remove later slides.
$sorted = $toDelete | Sort-Object -Descending

foreach ($idx in $sorted) {
    $p.Slides.Item($idx).Delete()
}

Write-Output ("Slides remaining: " + $p.Slides.Count)
